# EE is actually Estonian (now renamed to ET).
# The "EE / Ewe (Niger-Congo)" row was a mislabeled duplicate of Estonian
# (which already exists as the ET row). Remove the erroneous row entirely
# and carry over its "Detector" flag to the real Estonian (ET) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Estonian (ET) row already exists at row 16, but it was missing the
# "Detector" (column F) checkmark that belonged to it.
$ws.Range("F16").Value = "x"

# Remove the erroneous "EE / Ewe (Niger-Congo)" row (row 29). This shifts
# every row below it (TH/Thai, and the trailing blank styled row) up by one.
$ws.Rows.Item(29).Delete()

# Reflect the resulting selection position in the sheet view.
$ws.Range("F18").Select()
